# Commit: "#5: property aircraft done"
# The 建物 (Building) worksheet's "property_category" column (I) incorrectly
# held the value "land" for every data row. This fixes it to "building".

$wb = $excel.ActiveWorkbook

# The "建物" (Building) worksheet is sheet2.xml in the package
# (workbook.xml sheet order: 土地, 建物, 存款, 股票, 債務 -> index 2).
$ws = $wb.Worksheets.Item("建物")

# Data rows are 2 through 35 (row 1 is the header row).
$range = $ws.Range("I2:I35")
$range.Value = "building"
